$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number format on cells whose new values look numeric,
# so Excel stores them as text (matching the original inlineStr cell type)
# instead of auto-converting to a Number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Now write the actual values.
$ws.Range("D2").Value = "35.420.49"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.893.30"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").Value = "247.06"
$ws.Range("E5").Value = "  -3.04%  "
$ws.Range("D6").Value = "0.691"
$ws.Range("E6").Value = "  -4.12%  "
$ws.Range("E7").Value = "  -0.80%  "
$ws.Range("D8").Value = "43.98"
$ws.Range("E8").Value = "  +8.33%  "
$ws.Range("D9").Value = "0.350"
$ws.Range("E9").Value = "  -4.69%  "
$ws.Range("E10").Value = "  -4.11%  "
$ws.Range("D11").Value = "0.0972"
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("D12").Value = "13.15"
$ws.Range("E12").Value = "  +3.10%  "
$ws.Range("D13").Value = "2.171.29"
$ws.Range("D14").Value = "0.725"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").Value = "4.91"
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("D16").Value = "1.889.05"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("D17").Value = "35.425.68"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "73.27"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").Value = "0.0₃0821"
$ws.Range("E19").Value = "  -3.97%  "
$ws.Range("D20").Value = "246.63"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").Value = "12.82"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("E22").Value = "  -2.95%  "
$ws.Range("E24").Value = "  +6.33%  "
$ws.Range("E25").Value = "  -10.99%  "
$ws.Range("D26").Value = "165.23"
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("D27").Value = "8.47"
$ws.Range("E27").Value = "  -2.43%  "
$ws.Range("D28").Value = "18.35"
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("D29").Value = "0.127"
$ws.Range("E29").Value = "  -4.23%  "
$ws.Range("D30").Value = "4.128.40"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "1.80"
$ws.Range("E31").Value = "  +9.03%  "
$ws.Range("D32").Value = "4.25"
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("D33").Value = "0.0580"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("D34").Value = "4.25"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("D36").Value = "0.849"
$ws.Range("E36").Value = "  -6.70%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "1.60"
$ws.Range("E37").Value = "  -20.14%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "2.01"
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("D39").Value = "17.22"
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("D40").Value = "97.64"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").Value = "0.0675"
$ws.Range("E41").Value = "  +3.98%  "
$ws.Range("D42").Value = "0.0212"
$ws.Range("E42").Value = "  -2.14%  "
$ws.Range("D43").Value = "1.09"
$ws.Range("E43").Value = "  -2.19%  "
$ws.Range("D44").Value = "1.287.55"
$ws.Range("E44").Value = "  -3.63%  "
$ws.Range("D45").Value = "2.35"
$ws.Range("E45").Value = "  -2.57%  "
$ws.Range("D46").Value = "0.0806"
$ws.Range("E46").Value = "  +8.03%  "
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").Value = "12.06"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("E50").Value = "  -5.18%  "
$ws.Range("D51").Value = "43.09"
$ws.Range("E51").Value = "  -5.26%  "
